$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.977.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.263.74'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.87'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.652'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.69'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.452'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +6.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0983'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.74'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.58'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.600.48'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.69'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.839'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.265.89'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.918.60'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0988'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.89'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.12'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.32'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.39'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.05%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +21.49%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.61'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +9.83%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.91'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.45'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.96'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0686'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.97'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.29%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.49'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.89%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.20%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.20%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.35'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '98.61'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0951'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.39'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.86%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.37'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.449.35'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -9.48%  '
